# Auto-generated edit script
# Applies updated odds values to Sheet1 of the FlashScore workbook,
# matching the target XML diff (cell-by-cell Value2 assignment).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value2 = 1.3  # G2
$ws.Cells.Item(2, 9).Value2 = 9  # I2
$ws.Cells.Item(2, 12).Value2 = 7  # L2
$ws.Cells.Item(2, 15).Value2 = 1.1  # O2
$ws.Cells.Item(2, 16).Value2 = 7  # P2
$ws.Cells.Item(2, 28).Value2 = 9.5  # AB2
$ws.Cells.Item(2, 34).Value2 = 13  # AH2
$ws.Cells.Item(2, 40).Value2 = 29  # AN2
$ws.Cells.Item(2, 41).Value2 = 126  # AO2
$ws.Cells.Item(2, 43).Value2 = 51  # AQ2

# Row 3
$ws.Cells.Item(3, 8).Value2 = 3.5  # H3
$ws.Cells.Item(3, 9).Value2 = 1.9  # I3
$ws.Cells.Item(3, 17).Value2 = 1.93  # Q3
$ws.Cells.Item(3, 18).Value2 = 1.97  # R3
$ws.Cells.Item(3, 25).Value2 = 1.75  # Y3
$ws.Cells.Item(3, 26).Value2 = 2  # Z3
$ws.Cells.Item(3, 32).Value2 = 41  # AF3

# Row 4
$ws.Cells.Item(4, 17).Value2 = 1.67  # Q4
$ws.Cells.Item(4, 18).Value2 = 2.2  # R4

# Row 6
$ws.Cells.Item(6, 7).Value2 = 2.38  # G6
$ws.Cells.Item(6, 8).Value2 = 3.2  # H6
$ws.Cells.Item(6, 9).Value2 = 3.1  # I6
$ws.Cells.Item(6, 10).Value2 = 3.1  # J6
$ws.Cells.Item(6, 11).Value2 = 2.05  # K6
$ws.Cells.Item(6, 12).Value2 = 3.75  # L6
$ws.Cells.Item(6, 13).Value2 = 1.07  # M6
$ws.Cells.Item(6, 14).Value2 = 9  # N6
$ws.Cells.Item(6, 15).Value2 = 1.33  # O6
$ws.Cells.Item(6, 16).Value2 = 3.25  # P6
$ws.Cells.Item(6, 17).Value2 = 2.1  # Q6
$ws.Cells.Item(6, 18).Value2 = 1.7  # R6
$ws.Cells.Item(6, 19).Value2 = 3.1  # S6
$ws.Cells.Item(6, 20).Value2 = 1.37  # T6
$ws.Cells.Item(6, 21).Value2 = 3.75  # U6
$ws.Cells.Item(6, 22).Value2 = 1.25  # V6
$ws.Cells.Item(6, 23).Value2 = 1.44  # W6
$ws.Cells.Item(6, 24).Value2 = 2.63  # X6
$ws.Cells.Item(6, 25).Value2 = 1.8  # Y6
$ws.Cells.Item(6, 26).Value2 = 1.95  # Z6
$ws.Cells.Item(6, 28).Value2 = 11  # AB6
$ws.Cells.Item(6, 29).Value2 = 9.5  # AC6
$ws.Cells.Item(6, 30).Value2 = 23  # AD6
$ws.Cells.Item(6, 31).Value2 = 21  # AE6
$ws.Cells.Item(6, 32).Value2 = 29  # AF6
$ws.Cells.Item(6, 33).Value2 = 8.5  # AG6
$ws.Cells.Item(6, 37).Value2 = 301  # AK6
$ws.Cells.Item(6, 38).Value2 = 9  # AL6
$ws.Cells.Item(6, 39).Value2 = 15  # AM6
$ws.Cells.Item(6, 41).Value2 = 34  # AO6
$ws.Cells.Item(6, 43).Value2 = 34  # AQ6
$ws.Cells.Item(6, 44).Value2 = 1.6  # AR6
$ws.Cells.Item(6, 45).Value2 = 2.35  # AS6

# Row 7
$ws.Cells.Item(7, 13).Value2 = 1.14  # M7
$ws.Cells.Item(7, 14).Value2 = 5.5  # N7
$ws.Cells.Item(7, 30).Value2 = 21  # AD7
$ws.Cells.Item(7, 31).Value2 = 23  # AE7
$ws.Cells.Item(7, 38).Value2 = 7.5  # AL7
$ws.Cells.Item(7, 40).Value2 = 15  # AN7

# Row 12
$ws.Cells.Item(12, 7).Value2 = 1.09  # G12
$ws.Cells.Item(12, 8).Value2 = 8  # H12
$ws.Cells.Item(12, 9).Value2 = 32  # I12
$ws.Cells.Item(12, 10).Value2 = 1.38  # J12
$ws.Cells.Item(12, 11).Value2 = 3.05  # K12
$ws.Cells.Item(12, 12).Value2 = 18  # L12
$ws.Cells.Item(12, 13).Value2 = 1.02  # M12
$ws.Cells.Item(12, 14).Value2 = 11  # N12
$ws.Cells.Item(12, 15).Value2 = 1.12  # O12
$ws.Cells.Item(12, 16).Value2 = 5.6  # P12
$ws.Cells.Item(12, 17).Value2 = 1.38  # Q12
$ws.Cells.Item(12, 18).Value2 = 2.87  # R12
$ws.Cells.Item(12, 21).Value2 = 1.93  # U12
$ws.Cells.Item(12, 22).Value2 = 1.8  # V12
$ws.Cells.Item(12, 23).Value2 = 1.23  # W12
$ws.Cells.Item(12, 24).Value2 = 3.85  # X12
$ws.Cells.Item(12, 25).Value2 = 2.57  # Y12
$ws.Cells.Item(12, 26).Value2 = 1.45  # Z12
$ws.Cells.Item(12, 27).Value2 = 7  # AA12
$ws.Cells.Item(12, 28).Value2 = 6.3  # AB12
$ws.Cells.Item(12, 29).Value2 = 13.5  # AC12
$ws.Cells.Item(12, 30).Value2 = 6.1  # AD12
$ws.Cells.Item(12, 32).Value2 = 50  # AF12
$ws.Cells.Item(12, 33).Value2 = 11  # AG12
$ws.Cells.Item(12, 34).Value2 = 21  # AH12
$ws.Cells.Item(12, 36).Value2 = 300  # AJ12
$ws.Cells.Item(12, 38).Value2 = 55  # AL12
$ws.Cells.Item(12, 39).Value2 = 600  # AM12
$ws.Cells.Item(12, 40).Value2 = 150  # AN12
$ws.Cells.Item(12, 42).Value2 = 1000  # AP12
$ws.Cells.Item(12, 43).Value2 = 500  # AQ12

# Row 14
$ws.Cells.Item(14, 7).Value2 = 4.33  # G14
$ws.Cells.Item(14, 8).Value2 = 4.1  # H14
$ws.Cells.Item(14, 9).Value2 = 1.65  # I14
$ws.Cells.Item(14, 12).Value2 = 2.25  # L14
$ws.Cells.Item(14, 13).Value2 = 1.04  # M14
$ws.Cells.Item(14, 14).Value2 = 13  # N14
$ws.Cells.Item(14, 17).Value2 = 1.75  # Q14
$ws.Cells.Item(14, 18).Value2 = 2.05  # R14
$ws.Cells.Item(14, 23).Value2 = 1.33  # W14
$ws.Cells.Item(14, 24).Value2 = 3.25  # X14
$ws.Cells.Item(14, 25).Value2 = 1.75  # Y14
$ws.Cells.Item(14, 26).Value2 = 2  # Z14
$ws.Cells.Item(14, 37).Value2 = 201  # AK14
$ws.Cells.Item(14, 38).Value2 = 8  # AL14

# Row 15
$ws.Cells.Item(15, 17).Value2 = 2  # Q15
$ws.Cells.Item(15, 21).Value2 = 3.75  # U15
$ws.Cells.Item(15, 22).Value2 = 1.29  # V15
$ws.Cells.Item(15, 29).Value2 = 10  # AC15

# Row 16
$ws.Cells.Item(16, 7).Value2 = 1.48  # G16
$ws.Cells.Item(16, 8).Value2 = 4.1  # H16
$ws.Cells.Item(16, 9).Value2 = 7  # I16
$ws.Cells.Item(16, 10).Value2 = 2.1  # J16
$ws.Cells.Item(16, 11).Value2 = 2.1  # K16
$ws.Cells.Item(16, 12).Value2 = 8  # L16
$ws.Cells.Item(16, 13).Value2 = 1.06  # M16
$ws.Cells.Item(16, 14).Value2 = 10  # N16
$ws.Cells.Item(16, 15).Value2 = 1.36  # O16
$ws.Cells.Item(16, 16).Value2 = 3  # P16
$ws.Cells.Item(16, 17).Value2 = 2.2  # Q16
$ws.Cells.Item(16, 18).Value2 = 1.65  # R16
$ws.Cells.Item(16, 19).Value2 = 3.45  # S16
$ws.Cells.Item(16, 20).Value2 = 1.32  # T16
$ws.Cells.Item(16, 21).Value2 = 4  # U16
$ws.Cells.Item(16, 22).Value2 = 1.22  # V16
$ws.Cells.Item(16, 25).Value2 = 2.5  # Y16
$ws.Cells.Item(16, 26).Value2 = 1.5  # Z16
$ws.Cells.Item(16, 30).Value2 = 9.5  # AD16
$ws.Cells.Item(16, 34).Value2 = 8  # AH16
$ws.Cells.Item(16, 35).Value2 = 26  # AI16
$ws.Cells.Item(16, 36).Value2 = 101  # AJ16
$ws.Cells.Item(16, 38).Value2 = 13  # AL16
$ws.Cells.Item(16, 39).Value2 = 34  # AM16
$ws.Cells.Item(16, 40).Value2 = 23  # AN16
$ws.Cells.Item(16, 41).Value2 = 101  # AO16
$ws.Cells.Item(16, 42).Value2 = 67  # AP16
$ws.Cells.Item(16, 43).Value2 = 67  # AQ16

# Row 17
$ws.Cells.Item(17, 7).Value2 = 2.45  # G17
$ws.Cells.Item(17, 9).Value2 = 2.4  # I17
$ws.Cells.Item(17, 10).Value2 = 3.1  # J17
$ws.Cells.Item(17, 12).Value2 = 3  # L17
$ws.Cells.Item(17, 15).Value2 = 1.14  # O17
$ws.Cells.Item(17, 16).Value2 = 5.5  # P17
$ws.Cells.Item(17, 17).Value2 = 1.53  # Q17
$ws.Cells.Item(17, 18).Value2 = 2.4  # R17
$ws.Cells.Item(17, 19).Value2 = 1.83  # S17
$ws.Cells.Item(17, 20).Value2 = 1.98  # T17
$ws.Cells.Item(17, 21).Value2 = 2.25  # U17
$ws.Cells.Item(17, 22).Value2 = 1.57  # V17
$ws.Cells.Item(17, 28).Value2 = 15  # AB17
$ws.Cells.Item(17, 29).Value2 = 10  # AC17
$ws.Cells.Item(17, 30).Value2 = 26  # AD17
$ws.Cells.Item(17, 32).Value2 = 21  # AF17
$ws.Cells.Item(17, 34).Value2 = 7.5  # AH17
$ws.Cells.Item(17, 41).Value2 = 26  # AO17
$ws.Cells.Item(17, 42).Value2 = 19  # AP17

# Row 20
$ws.Cells.Item(20, 7).Value2 = 1.7  # G20
$ws.Cells.Item(20, 8).Value2 = 3.9  # H20
$ws.Cells.Item(20, 9).Value2 = 4.33  # I20
$ws.Cells.Item(20, 10).Value2 = 2.25  # J20
$ws.Cells.Item(20, 12).Value2 = 4.33  # L20
$ws.Cells.Item(20, 25).Value2 = 1.5  # Y20
$ws.Cells.Item(20, 26).Value2 = 2.5  # Z20
$ws.Cells.Item(20, 27).Value2 = 11  # AA20
$ws.Cells.Item(20, 28).Value2 = 11  # AB20
$ws.Cells.Item(20, 30).Value2 = 15  # AD20
$ws.Cells.Item(20, 33).Value2 = 19  # AG20
$ws.Cells.Item(20, 39).Value2 = 26  # AM20
$ws.Cells.Item(20, 41).Value2 = 41  # AO20
$ws.Cells.Item(20, 42).Value2 = 29  # AP20
$ws.Cells.Item(20, 43).Value2 = 29  # AQ20

# Row 21
$ws.Cells.Item(21, 7).Value2 = 2.25  # G21
$ws.Cells.Item(21, 8).Value2 = 3.5  # H21
$ws.Cells.Item(21, 10).Value2 = 2.88  # J21
$ws.Cells.Item(21, 11).Value2 = 2.25  # K21
$ws.Cells.Item(21, 12).Value2 = 3.5  # L21
$ws.Cells.Item(21, 13).Value2 = 1.04  # M21
$ws.Cells.Item(21, 14).Value2 = 13  # N21
$ws.Cells.Item(21, 15).Value2 = 1.22  # O21
$ws.Cells.Item(21, 16).Value2 = 4  # P21
$ws.Cells.Item(21, 17).Value2 = 1.75  # Q21
$ws.Cells.Item(21, 18).Value2 = 2.05  # R21
$ws.Cells.Item(21, 21).Value2 = 2.75  # U21
$ws.Cells.Item(21, 22).Value2 = 1.4  # V21
$ws.Cells.Item(21, 23).Value2 = 1.33  # W21
$ws.Cells.Item(21, 24).Value2 = 3.25  # X21
$ws.Cells.Item(21, 25).Value2 = 1.62  # Y21
$ws.Cells.Item(21, 26).Value2 = 2.2  # Z21
$ws.Cells.Item(21, 27).Value2 = 9.5  # AA21
$ws.Cells.Item(21, 28).Value2 = 12  # AB21
$ws.Cells.Item(21, 31).Value2 = 17  # AE21
$ws.Cells.Item(21, 32).Value2 = 23  # AF21
$ws.Cells.Item(21, 33).Value2 = 13  # AG21
$ws.Cells.Item(21, 34).Value2 = 7  # AH21
$ws.Cells.Item(21, 35).Value2 = 13  # AI21
$ws.Cells.Item(21, 36).Value2 = 41  # AJ21
$ws.Cells.Item(21, 37).Value2 = 151  # AK21
$ws.Cells.Item(21, 38).Value2 = 11  # AL21
$ws.Cells.Item(21, 39).Value2 = 17  # AM21
$ws.Cells.Item(21, 41).Value2 = 29  # AO21
$ws.Cells.Item(21, 43).Value2 = 29  # AQ21

# Row 22
$ws.Cells.Item(22, 11).Value2 = 2.5  # K22
$ws.Cells.Item(22, 13).Value2 = 1.02  # M22
$ws.Cells.Item(22, 14).Value2 = 19  # N22
$ws.Cells.Item(22, 17).Value2 = 1.5  # Q22
$ws.Cells.Item(22, 18).Value2 = 2.5  # R22
$ws.Cells.Item(22, 19).Value2 = 1.8  # S22
$ws.Cells.Item(22, 20).Value2 = 2.05  # T22
$ws.Cells.Item(22, 21).Value2 = 2.2  # U22
$ws.Cells.Item(22, 22).Value2 = 1.62  # V22
$ws.Cells.Item(22, 23).Value2 = 1.25  # W22
$ws.Cells.Item(22, 24).Value2 = 3.75  # X22
$ws.Cells.Item(22, 28).Value2 = 12  # AB22
$ws.Cells.Item(22, 38).Value2 = 17  # AL22
$ws.Cells.Item(22, 39).Value2 = 23  # AM22

# Row 23
$ws.Cells.Item(23, 7).Value2 = 3.1  # G23
$ws.Cells.Item(23, 9).Value2 = 2.15  # I23
$ws.Cells.Item(23, 29).Value2 = 12  # AC23

# Row 25
$ws.Cells.Item(25, 7).Value2 = 3.75  # G25
$ws.Cells.Item(25, 8).Value2 = 3.4  # H25
$ws.Cells.Item(25, 9).Value2 = 2  # I25
$ws.Cells.Item(25, 10).Value2 = 4.33  # J25
$ws.Cells.Item(25, 11).Value2 = 2.1  # K25
$ws.Cells.Item(25, 12).Value2 = 2.75  # L25
$ws.Cells.Item(25, 13).Value2 = 1.06  # M25
$ws.Cells.Item(25, 14).Value2 = 10  # N25
$ws.Cells.Item(25, 17).Value2 = 2.08  # Q25
$ws.Cells.Item(25, 18).Value2 = 1.73  # R25
$ws.Cells.Item(25, 27).Value2 = 10  # AA25
$ws.Cells.Item(25, 28).Value2 = 19  # AB25
$ws.Cells.Item(25, 31).Value2 = 34  # AE25
$ws.Cells.Item(25, 33).Value2 = 9  # AG25
$ws.Cells.Item(25, 39).Value2 = 9  # AM25
$ws.Cells.Item(25, 41).Value2 = 17  # AO25
$ws.Cells.Item(25, 42).Value2 = 17  # AP25

# Row 27
$ws.Cells.Item(27, 7).Value2 = 1.72  # G27
$ws.Cells.Item(27, 8).Value2 = 3.85  # H27
$ws.Cells.Item(27, 9).Value2 = 4.2  # I27
$ws.Cells.Item(27, 10).Value2 = 2.25  # J27
$ws.Cells.Item(27, 11).Value2 = 2.27  # K27
$ws.Cells.Item(27, 12).Value2 = 4.4  # L27
$ws.Cells.Item(27, 13).Value2 = 1.04  # M27
$ws.Cells.Item(27, 14).Value2 = 8.5  # N27
$ws.Cells.Item(27, 15).Value2 = 1.22  # O27
$ws.Cells.Item(27, 16).Value2 = 3.9  # P27
$ws.Cells.Item(27, 17).Value2 = 1.65  # Q27
$ws.Cells.Item(27, 18).Value2 = 2.1  # R27
$ws.Cells.Item(27, 21).Value2 = 2.57  # U27
$ws.Cells.Item(27, 22).Value2 = 1.45  # V27
$ws.Cells.Item(27, 23).Value2 = 1.34  # W27
$ws.Cells.Item(27, 24).Value2 = 3  # X27
$ws.Cells.Item(27, 25).Value2 = 1.65  # Y27
$ws.Cells.Item(27, 26).Value2 = 2.1  # Z27
$ws.Cells.Item(27, 27).Value2 = 8.5  # AA27
$ws.Cells.Item(27, 30).Value2 = 14  # AD27
$ws.Cells.Item(27, 31).Value2 = 12.5  # AE27
$ws.Cells.Item(27, 32).Value2 = 22  # AF27
$ws.Cells.Item(27, 33).Value2 = 8.5  # AG27
$ws.Cells.Item(27, 34).Value2 = 7.6  # AH27
$ws.Cells.Item(27, 35).Value2 = 14  # AI27
$ws.Cells.Item(27, 37).Value2 = 350  # AK27
$ws.Cells.Item(27, 38).Value2 = 14  # AL27
$ws.Cells.Item(27, 39).Value2 = 25  # AM27
$ws.Cells.Item(27, 40).Value2 = 13.5  # AN27
$ws.Cells.Item(27, 41).Value2 = 65  # AO27
$ws.Cells.Item(27, 43).Value2 = 37  # AQ27

# Row 29
$ws.Cells.Item(29, 15).Value2 = 1.11  # O29
$ws.Cells.Item(29, 16).Value2 = 6.5  # P29
$ws.Cells.Item(29, 17).Value2 = 1.4  # Q29
$ws.Cells.Item(29, 18).Value2 = 2.88  # R29
$ws.Cells.Item(29, 21).Value2 = 1.91  # U29
$ws.Cells.Item(29, 22).Value2 = 1.8  # V29
